# edit.ps1 - apply week05 response-document edits:
#   1) Rewrite the "Solution: Create FOR loop ..." paragraph (union solution)
#      into six separate runs:
#        "Solution:" / " " / "Add set1 & set2 to set3." / " " /
#        "The set will automatically throw out duplicates." / " "
#   2) Fill the blank paragraph right after the Find-Pairs interview
#      question with a new sentence.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: locate a paragraph index whose text matches a wildcard pattern
# ---------------------------------------------------------------------
function Find-ParagraphIndex([string]$pattern) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -like $pattern) {
            return $i
        }
    }
    return -1
}

# =======================================================================
# Change 1: "Solution: Create FOR loop ..." -> split into multiple runs
# =======================================================================

$solutionIdx = Find-ParagraphIndex("*Create FOR loop to add values from sets 1 & 2 to set3*")

if ($solutionIdx -ne -1) {
    $solutionPara = $d.Paragraphs($solutionIdx)

    # Completely empty the paragraph (remove all its text, but keep the
    # trailing paragraph mark / paragraph-level identity).
    $fullRange = $d.Range($solutionPara.Range.Start, $solutionPara.Range.End - 1)
    $fullRange.Delete()

    # All pieces that must end up as separate <w:r> runs, in order. The
    # very last piece is typed straight into the original (now emptied)
    # paragraph so that paragraph keeps its identity / paragraph
    # properties; the earlier pieces are built in brand-new paragraphs
    # placed *before* it and then folded forward into it one
    # paragraph-mark at a time (each fold keeps the *later* paragraph's
    # mark, i.e. the original one, alive).
    $piecesBeforeLast = @("Solution:", " ", "Add set1 & set2 to set3.", " ", "The set will automatically throw out duplicates.")
    $lastPiece = " "

    $insertAfterIdx = $solutionIdx - 1
    foreach ($piece in $piecesBeforeLast) {
        $srcPara = $d.Paragraphs($insertAfterIdx)
        $srcPara.Range.InsertParagraphAfter()
        $insertAfterIdx = $insertAfterIdx + 1
        $d.Paragraphs($insertAfterIdx).Range.InsertAfter($piece)
    }

    # The original (now empty) paragraph has been pushed down to index
    # ($insertAfterIdx + 1); put the final piece straight into it.
    $lastIdx = $insertAfterIdx + 1
    $d.Paragraphs($lastIdx).Range.InsertAfter($lastPiece)

    # Fold every inserted paragraph forward into the original one by
    # repeatedly deleting the paragraph mark that precedes it - this
    # merges content while letting the later (original) paragraph's
    # mark - and thus its paragraph-level identity - survive.
    while ($lastIdx -gt $solutionIdx) {
        $prevIdx = $lastIdx - 1
        $prevPara = $d.Paragraphs($prevIdx)
        $mark = $d.Range($prevPara.Range.End - 1, $prevPara.Range.End)
        $mark.Delete()
        $lastIdx = $prevIdx
    }
}

# =======================================================================
# Change 2: fill the blank paragraph after the Find-Pairs question
# =======================================================================

$questionIdx = Find-ParagraphIndex("*Find Pairs*30 seconds*")

if ($questionIdx -ne -1) {
    $blankIdx = $questionIdx + 1
    $blankPara = $d.Paragraphs($blankIdx)
    if ($blankPara.Range.Text.Trim() -eq "") {
        $blankPara.Range.InsertAfter("Check to see if the other values are reversed versions of the current one.")
    }
}
